# Update countries & provincias Spain
# - Reorders two country rows (Costa de Marfil, Burkina Faso move earlier in the
#   country list) and refreshes their case counts, plus refreshed totals for a
#   handful of other countries, and bumps the "last updated" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 23:52"

# Updated totals for existing countries (no row movement)
$ws.Range("B4").Value = 529154
$ws.Range("C4").Value = 26278
$ws.Range("E4").Value = 479252
$ws.Range("G4").Value = 1713
$ws.Range("H4").Value = 20460

$ws.Range("B5").Value = 163027
$ws.Range("C5").Value = 4754
$ws.Range("E5").Value = 87312
$ws.Range("G5").Value = 525
$ws.Range("H5").Value = 16606

$ws.Range("B16").Value = 23318
$ws.Range("C16").Value = 1170
$ws.Range("E16").Value = 16356
$ws.Range("G16").Value = 84
$ws.Range("H16").Value = 653

$ws.Range("B29").Value = 6409
$ws.Range("C29").Value = 95
$ws.Range("E29").Value = 6258

# Rows 91-95: Costa de Marfil and Burkina Faso move earlier in the country
# ranking (inserted right after Oman, and right after Niger respectively),
# pushing Uruguay, Niger and Banglades down one slot each, and all five rows
# get refreshed case counts.
$ws.Range("A91").Value = "Costa de Marfil"
$ws.Range("B91").Value = 533
$ws.Range("C91").Value = 53
$ws.Range("D91").Value = 58
$ws.Range("E91").Value = 471
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 4

$ws.Range("A92").Value = "Uruguay"
$ws.Range("B92").Value = 494
$ws.Range("C92").Value = 21
$ws.Range("D92").Value = 214
$ws.Range("E92").Value = 273
$ws.Range("F92").Value = 15
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 7

$ws.Range("A93").Value = "Niger"
$ws.Range("B93").Value = 491
$ws.Range("C93").Value = 53
$ws.Range("D93").Value = 41
$ws.Range("E93").Value = 439
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 11

$ws.Range("A94").Value = "Burkina Faso"
$ws.Range("B94").Value = 484
$ws.Range("C94").Value = 36
$ws.Range("D94").Value = 155
$ws.Range("E94").Value = 302
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 27

$ws.Range("A95").Value = "Banglades"
$ws.Range("B95").Value = 482
$ws.Range("C95").Value = 58
$ws.Range("D95").Value = 36
$ws.Range("E95").Value = 416
$ws.Range("F95").Value = 1
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 30
